# Applies new applicant rows 203-224 to the admissions worksheet (qabul.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns F, I, J, K hold numeric-looking / date-looking text (JSHIR, phone
# numbers, dates) that must stay plain text, matching the rest of the sheet.
$textForcedColumns = @("F", "I", "J", "K")

$newApplicants = @(
    @{ Row = 203; A = "Rasuljonov Fazliddin"; B = "Bugalteriya hisobi"; C = "O'zbek tili"; D = "Kunduzgi"; E = "AD2481404"; F = "50401076530012"; G = "Toshkent shahri"; H = "Mirobod tumani"; I = "998942433422"; J = "+998942433422"; K = "2025-07-15" },
    @{ Row = 204; A = "Bazarbayev Samir Quanishbay o'g'li"; B = "Yurisprudensiya"; C = "O'zbek tili"; D = "Kunduzgi"; E = "AD6256423"; F = "50611077350054"; G = "Qoraqalpogʻiston Respublikasi"; H = "Beruniy tumani"; I = "998997295701"; J = "+998997295701"; K = "2025-07-16" },
    @{ Row = 205; A = "Avaznazarov Diyorbek"; B = "Yurisprudensiya"; C = "O'zbek tili"; D = "Kunduzgi"; E = "AD3275498"; F = "52910065710015"; G = "Qashqadaryo viloyati"; H = "Mirishkor tumani"; I = "998908784346"; J = "+998950297101"; K = "2025-07-16" },
    @{ Row = 206; A = "Abdusattorova Shahloxon"; B = "Yurisprudensiya"; C = "O'zbek tili"; D = "Kunduzgi"; E = "AD0332937"; F = "62808046940058"; G = "Fargona viloyati"; H = "Quva tumani"; I = "998955805054"; J = "+998916655156"; K = "2025-07-16" },
    @{ Row = 207; A = "Xasanboyeva Marjona Asatilla qizi"; B = "Hayot faoliyati xavfsizligi"; C = "O'zbek tili"; D = "Kunduzgi"; E = "AD6725608"; F = "62702086580013"; G = "Toshkent shahri"; H = "Yashnaobod tumani"; I = "998993019919"; J = "+998933519919"; K = "2025-07-16" },
    @{ Row = 208; A = "Sunnatboyev Asilbek Lutfulloyevich"; B = "Yurisprudensiya"; C = "Rus tili"; D = "Kunduzgi"; E = "AD6172894"; F = "51102058540026"; G = "Navoiy viloyati"; H = "Navoiy shahri"; I = "998930022274"; J = "+998931512274"; K = "2025-07-16" },
    @{ Row = 209; A = "Rahmonov Humoyun Xayitmurodovich"; B = "Yurisprudensiya"; C = "O'zbek tili"; D = "Kunduzgi"; E = "AB2339324"; F = "30808985700024"; G = "Qashqadaryo viloyati"; H = "Kasbi tumani"; I = "998882090096"; J = "+998881113261"; K = "2025-07-16" },
    @{ Row = 210; A = "Sobirjonov Saidjon Obidjon o'g'li"; B = "Yurisprudensiya"; C = "Rus tili"; D = "Kunduzgi"; E = "AC2641623"; F = "51101045950016"; G = "Namangan viloyati"; H = "Yangiqoʻrgʻon tumani"; I = "998933772113"; J = "+998772341110"; K = "2025-07-16" },
    @{ Row = 211; A = "Jasmin Batirova"; B = "Yurisprudensiya"; C = "O'zbek tili"; D = "Kunduzgi"; E = "AD4597593"; F = "61802077190089"; G = "Toshkent shahri"; H = "Yunusobod tumani"; I = "998886872777"; J = "+998886872777"; K = "2025-07-16" },
    @{ Row = 212; A = "Maxmadmurodov Karim Orifjon ogli"; B = "Ijtimoiy ish"; C = "O'zbek tili"; D = "Kunduzgi"; E = "AD3677109"; F = "50203076580028"; G = "Toshkent shahri"; H = "Yashnaobod tumani"; I = "998997668641"; J = "+998997668699"; K = "2025-07-16" },
    @{ Row = 213; A = "Shohjahon Botirov"; B = "Yurisprudensiya"; C = "O'zbek tili"; D = "Kunduzgi"; E = "AD6006127"; F = "51107087140013"; G = "Toshkent shahri"; H = "Yunusobod tumani"; I = "998950261107"; J = "+998909451411"; K = "2025-07-16" },
    @{ Row = 214; A = "Adashboyeva Fotima Qodirali qizi"; B = "Yurisprudensiya"; C = "O'zbek tili"; D = "Kunduzgi"; E = "AD1220503"; F = "62508005990011"; G = "Namangan viloyati"; H = "Mingbuloq tumani"; I = "79126723290"; J = "+998333031325"; K = "2025-07-16" },
    @{ Row = 215; A = "Dilshodova Jasmina Dilshodovna"; B = "Yurisprudensiya"; C = "Rus tili"; D = "Kunduzgi"; E = "AD2836503"; F = "61709066500066"; G = "Toshkent shahri"; H = "Uchtepa tumani"; I = "998773788878"; J = "+998333704271"; K = "2025-07-17" },
    @{ Row = 216; A = "Murodova Orzigul Xolmumin qizi"; B = "Yurisprudensiya"; C = "Rus tili"; D = "Kunduzgi"; E = "AD4184431"; F = "62207076080059"; G = "Toshkent shahri"; H = "Chilonzor tumani"; I = "998931312207"; J = "+998931312207"; K = "2025-07-17" },
    @{ Row = 217; A = "Shamsiddin Bahodirov Anvar ogli"; B = "Yurisprudensiya"; C = "O'zbek tili"; D = "Kunduzgi"; E = "AD1709792"; F = "52201056430028"; G = "Sirdaryo viloyati"; H = "Mirzaobod tumani"; I = "998994602082"; J = "+998990372230"; K = "2025-07-17" },
    @{ Row = 218; A = "Abdiyev Bunyod"; B = "Yurisprudensiya"; C = "O'zbek tili"; D = "Kunduzgi"; E = "AD8597324"; F = "50604086450020"; G = "Sirdaryo viloyati"; H = "Guliston tumani"; I = "998999090257"; J = "+998911010257"; K = "2025-07-17" },
    @{ Row = 219; A = "Xusanboyev Murodbek Sherbek o'g'li"; B = "Yurisprudensiya"; C = "O'zbek tili"; D = "Kunduzgi"; E = "AD9591406"; F = "52709075530012"; G = "Jizzax viloyati"; H = "Yangiobod tumani"; I = "+998996603300"; J = "+998995760701"; K = "2025-07-17" },
    @{ Row = 220; A = "Rasulova Sevinch Kamoliddin qizi"; B = "Yurisprudensiya"; C = "O'zbek tili"; D = "Kunduzgi"; E = "AD3631826"; F = "61206075530012"; G = "Jizzax viloyati"; H = "Yangiobod tumani"; I = "998972011207"; J = "+998972011207"; K = "2025-07-17" },
    @{ Row = 221; A = "Murodov Firdavs Ozodbek ogli"; B = "Yurisprudensiya"; C = "O'zbek tili"; D = "Kunduzgi"; E = "AD2765996"; F = "50907065530015"; G = "Jizzax viloyati"; H = "Yangiobod tumani"; I = "998904020066"; J = "+998904020066"; K = "2025-07-17" },
    @{ Row = 222; A = "Botirova Gulmira Quvondiq qizi"; B = "Yurisprudensiya"; C = "O'zbek tili"; D = "Kunduzgi"; E = "AB4552654"; F = "61702007100017"; G = "Xorazm viloyati"; H = "Bogʻot tumani"; I = "998942331799"; J = "+998932890124"; K = "2025-07-17" },
    @{ Row = 223; A = "Amriyeva Nozima"; B = "Yurisprudensiya"; C = "O'zbek tili"; D = "Kunduzgi"; E = "AD8022695"; F = "61807076110030"; G = "Samarqand viloyati"; H = "Samarqand shahri"; I = "998939671807"; J = "+998939671807"; K = "2025-07-18" },
    @{ Row = 224; A = "Latipov Akmal Akobirovich"; B = "Yurisprudensiya"; C = "O'zbek tili"; D = "Kunduzgi"; E = "AD0713461"; F = "32808871110038"; G = "Buxoro viloyati"; H = "Jondor tumani"; I = "+998993103666"; J = "+998993103666"; K = "2025-07-18" }
)

foreach ($applicant in $newApplicants) {
    $r = $applicant.Row
    foreach ($col in @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "K")) {
        $cell = $ws.Range("$col$r")
        if ($textForcedColumns -contains $col) {
            # Force text interpretation so numeric/date-like strings (e.g.
            # "998942433422" or "2025-07-15") are not auto-converted by Excel,
            # then clear the resulting format so no style is left on the cell.
            $cell.NumberFormat = "@"
            $cell.Value = $applicant[$col]
            $cell.ClearFormats()
        } else {
            $cell.Value = $applicant[$col]
        }
    }
}
